# The workbook's only sheet ("Template_NC") has a single header row (row 1)
# that originally spanned columns A:AN (40 columns). This edit:
#   1. Rewrites the header row down to 16 columns (A:P) with new header text
#      (re-purposing/renaming many of the old headers, e.g. nota contabila
#      automatic-numbering columns, customer balance/difference columns, etc.)
#   2. Removes all the now-unused trailing columns (Q:AN) from row 1.
#   3. Leaves the active selection on cell O2 (just under the new "GCI" header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values for columns A1 through P1, in order.
$headers = @(
    "BR",
    "Statutory_GL",
    "Journal",
    "Open_Item",
    "File_Ref",
    "Date",
    "Month",
    "Year",
    "Amount",
    "data_Description",
    "Post_Date",
    "TC",
    "Foreign Amount",
    "Foreign Currency",
    "GCI",
    "Customer Name"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Drop the old trailing headers (previously columns Q1:AN1) - no longer needed.
$ws.Range("Q1:AN1").ClearContents()

# Match the saved selection state from the edit.
$ws.Range("O2").Select()
